# Apply updated policy values to Sheet1.
# Column F = "Entertainment(Indoor)", Column H = "Industries"
#
# Rows 26-124: F changes from 0.6 -> 0.8
# Rows 125-132: F changes from 0   -> 0.2
# Rows 34-176: H changes from 1   -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 26; $r -le 124; $r++) {
    $ws.Range("F$r").Value = 0.8
}

for ($r = 125; $r -le 132; $r++) {
    $ws.Range("F$r").Value = 0.2
}

for ($r = 34; $r -le 176; $r++) {
    $ws.Range("H$r").Value = 0
}
